$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 290.77777
$ws.Range("I33").Value = 306.7647
$ws.Range("K33").Value = 306.7647
$ws.Range("M33").Value = -77.7647
$ws.Range("H42").Value = 771.7778
$ws.Range("I42").Value = 226.5
$ws.Range("K42").Value = 679.5
$ws.Range("M42").Value = -449.5
$ws.Range("H101").Value = 4772.8
$ws.Range("I101").Value = 3466.25
$ws.Range("K101").Value = 10398.75
$ws.Range("M101").Value = -8776.75
$ws.Range("H115").Value = 2638
$ws.Range("I115").Value = 1434.5
$ws.Range("J115").Value = 6248.5
$ws.Range("K115").Value = 4303.5
$ws.Range("L115").Value = 18745.5
$ws.Range("M115").Value = -2736.5
$ws.Range("N115").Value = -21879.5
$ws.Range("H131").Value = 1198
$ws.Range("I131").Value = 1198
$ws.Range("K131").Value = 3594
$ws.Range("M131").Value = 1446
$ws.Range("H138").Value = 3223.9321
$ws.Range("J138").Value = 5189.0625
$ws.Range("L138").Value = 15567.1875
$ws.Range("N138").Value = -25847.1875
$ws.Range("H141").Value = 2789.4468
$ws.Range("I141").Value = 2639.4092
$ws.Range("K141").Value = 7918.2276
$ws.Range("M141").Value = -2738.2276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2056.6
$ws.Range("I45").Value = 1838.2858
$ws.Range("J45").Value = 2566
$ws.Range("K45").Value = 1838.2858
$ws.Range("L45").Value = 2566
$ws.Range("M45").Value = -1461.2858
$ws.Range("N45").Value = -3320
$ws.Range("H88").Value = 10978.417
$ws.Range("I88").Value = 9500
$ws.Range("J88").Value = 11274.1
$ws.Range("K88").Value = 9500
$ws.Range("L88").Value = 11274.1
$ws.Range("M88").Value = -9094
$ws.Range("N88").Value = -12086.1
$ws.Range("H91").Value = 10978.417
$ws.Range("I91").Value = 9500
$ws.Range("J91").Value = 11274.1
$ws.Range("K91").Value = 9500
$ws.Range("L91").Value = 11274.1
$ws.Range("M91").Value = -8096
$ws.Range("N91").Value = -14082.1
$ws.Range("H97").Value = 845.55554
$ws.Range("I97").Value = 801.375
$ws.Range("J97").Value = 1199
$ws.Range("K97").Value = 801.375
$ws.Range("L97").Value = 1199
$ws.Range("M97").Value = -305.375
$ws.Range("N97").Value = -2191
$ws.Range("H122").Value = 6085.8887
$ws.Range("I122").Value = 4994.6
$ws.Range("J122").Value = 7450
$ws.Range("K122").Value = 14983.8
$ws.Range("L122").Value = 22350
$ws.Range("M122").Value = -12533.8
$ws.Range("N122").Value = -27250
$ws.Range("H132").Value = 1256894.9
$ws.Range("I132").Value = 1685429
$ws.Range("J132").Value = 154949.86
$ws.Range("K132").Value = 5056287
$ws.Range("L132").Value = 464849.58
$ws.Range("M132").Value = -5053757
$ws.Range("N132").Value = -469909.58

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 10662.6
$ws.Range("I86").Value = 9656.357
$ws.Range("K86").Value = 9656.357
$ws.Range("M86").Value = -8533.357
$ws.Range("H89").Value = 10662.6
$ws.Range("I89").Value = 9656.357
$ws.Range("K89").Value = 48281.785
$ws.Range("M89").Value = -42665.785

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3298.2307
$ws.Range("I122").Value = 2619.7778
$ws.Range("J122").Value = 4824.75
$ws.Range("K122").Value = 7859.3334
$ws.Range("L122").Value = 14474.25
$ws.Range("M122").Value = -5409.3334
$ws.Range("N122").Value = -19374.25
$ws.Range("H132").Value = 5814.027
$ws.Range("I132").Value = 4771.645
$ws.Range("J132").Value = 11199.667
$ws.Range("K132").Value = 14314.935
$ws.Range("L132").Value = 33599.001
$ws.Range("M132").Value = -11784.935
$ws.Range("N132").Value = -38659.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 6376.769
$ws.Range("I108").Value = 1079.6
$ws.Range("K108").Value = 3238.8
$ws.Range("M108").Value = -358.7999999999997
$ws.Range("H109").Value = 12557.692
$ws.Range("I109").Value = 9906.25
$ws.Range("K109").Value = 29718.75
$ws.Range("M109").Value = -28678.75
$ws.Range("H131").Value = 28890896
$ws.Range("I131").Value = 31373284
$ws.Range("J131").Value = 25644696
$ws.Range("K131").Value = 94119852
$ws.Range("L131").Value = 76934088
$ws.Range("M131").Value = -94114812
$ws.Range("N131").Value = -76944168

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 57666
$ws.Range("J111").Value = 57666
$ws.Range("L111").Value = 57666
$ws.Range("N111").Value = -63800
$ws.Range("H112").Value = 67000
$ws.Range("J112").Value = 67000
$ws.Range("L112").Value = 67000
$ws.Range("N112").Value = -69216
$ws.Range("H132").Value = 55557868
$ws.Range("I132").Value = 66668564
$ws.Range("J132").Value = 4402.6665
$ws.Range("K132").Value = 200005692
$ws.Range("L132").Value = 13207.9995
$ws.Range("M132").Value = -200003162
$ws.Range("N132").Value = -18267.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6053.3335
$ws.Range("I40").Value = 4838.8125
$ws.Range("J40").Value = 7819.909
$ws.Range("K40").Value = 4838.8125
$ws.Range("L40").Value = 7819.909
$ws.Range("M40").Value = -4702.8125
$ws.Range("N40").Value = -8091.909
$ws.Range("H68").Value = 2663.0527
$ws.Range("J68").Value = 3366.3333
$ws.Range("L68").Value = 3366.3333
$ws.Range("N68").Value = -4864.3333
$ws.Range("H71").Value = 2663.0527
$ws.Range("J71").Value = 3366.3333
$ws.Range("L71").Value = 16831.6665
$ws.Range("N71").Value = -24319.6665
$ws.Range("H106").Value = 30816.166
$ws.Range("J106").Value = 30816.166
$ws.Range("L106").Value = 30816.166
$ws.Range("N106").Value = -33340.166
$ws.Range("H122").Value = 3537
$ws.Range("I122").Value = 3482.6667
$ws.Range("J122").Value = 3700
$ws.Range("K122").Value = 10448.0001
$ws.Range("L122").Value = 11100
$ws.Range("M122").Value = -7998.000100000001
$ws.Range("N122").Value = -16000
$ws.Range("H132").Value = 3983.9
$ws.Range("I132").Value = 3404.1052
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 10212.3156
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -7682.3156
$ws.Range("N132").Value = -50060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12000
$ws.Range("J62").Value = 12874.833
$ws.Range("L62").Value = 12874.833
$ws.Range("N62").Value = -14122.833
$ws.Range("H65").Value = 12000
$ws.Range("J65").Value = 12874.833
$ws.Range("L65").Value = 64374.165
$ws.Range("N65").Value = -70614.16500000001
$ws.Range("H81").Value = 1392.8276
$ws.Range("I81").Value = 1331.68
$ws.Range("J81").Value = 1775
$ws.Range("K81").Value = 2663.36
$ws.Range("L81").Value = 3550
$ws.Range("M81").Value = -1602.36
$ws.Range("N81").Value = -5672
$ws.Range("H84").Value = 1392.8276
$ws.Range("I84").Value = 1331.68
$ws.Range("J84").Value = 1775
$ws.Range("K84").Value = 13316.8
$ws.Range("L84").Value = 17750
$ws.Range("M84").Value = -8012.800000000001
$ws.Range("N84").Value = -28358
$ws.Range("H104").Value = 56459.89
$ws.Range("J104").Value = 56459.89
$ws.Range("L104").Value = 56459.89
$ws.Range("N104").Value = -63447.89
$ws.Range("H107").Value = 486.9655
$ws.Range("I107").Value = 495.5
$ws.Range("J107").Value = 413
$ws.Range("K107").Value = 1486.5
$ws.Range("L107").Value = 1239
$ws.Range("M107").Value = 433.5
$ws.Range("N107").Value = -5079
$ws.Range("H122").Value = 1853.5834
$ws.Range("I122").Value = 1810.9375
$ws.Range("K122").Value = 5432.8125
$ws.Range("M122").Value = -2982.8125
$ws.Range("H132").Value = 7594.92
$ws.Range("I132").Value = 4494.55
$ws.Range("J132").Value = 19996.4
$ws.Range("K132").Value = 13483.65
$ws.Range("L132").Value = 59989.2
$ws.Range("M132").Value = -10953.65
$ws.Range("N132").Value = -65049.2
